$d = $word.ActiveDocument

# The error message raised when a bookmark block is missing its
# [ENDBOOKMARK] tag now also reports which bookmark was being parsed,
# e.g. "... while parsing m:bookmark 'bookmark1'".
$oldText = "Invalid block: Unexpected tag EOF missing [ENDBOOKMARK]"
$suffix = " while parsing m:bookmark 'bookmark1'"
$newText = $oldText + $suffix

# Guard against double-application: only touch the document if the
# longer (already-fixed) message is not present yet.
$already = $d.Content
$alreadyFound = $already.Find.Execute($newText, $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)

if (-not $alreadyFound) {
    $rng = $d.Content
    $found = $rng.Find.Execute(
        $oldText, $true, $false, $false, $false, $false,
        $true, 1, $false, "", 0)

    if ($found) {
        # Assign Range.Text directly (instead of Find.Execute's Replace
        # parameter) so straight apostrophes are preserved and Word's
        # AutoCorrect "smart quotes" feature does not turn them into
        # curly quotes.
        $rng.Text = $newText
    }
}
